# Fruta / hortaliza, semanal
#
# Inserts a new weekly price-report row for "Membrillo" (Feria Lagunitas de
# Puerto Montt) ahead of the existing row 47, pushing the subsequent rows
# (old rows 47-132) down by one (new rows 48-133). The new row reuses the
# same product/quality/origin metadata as the row that used to sit at 47,
# but carries its own date (2023-03-29) and volume (120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 47; everything below (old rows 47-132)
# shifts down to rows 48-133, carrying its formatting with it.
$ws.Rows(47).Insert()

# Populate the newly inserted row 47.
$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 45014
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = "Fruta"
$ws.Cells.Item(47, 7).Value = 100104
$ws.Cells.Item(47, 8).Value = "Frutos de pepita"
$ws.Cells.Item(47, 9).Value = 100104003
$ws.Cells.Item(47, 10).Value = "Membrillo"
$ws.Cells.Item(47, 11).Value = "Champion"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 120
$ws.Cells.Item(47, 14).Value = 17000
$ws.Cells.Item(47, 15).Value = 18000
$ws.Cells.Item(47, 16).Value = 17500
$ws.Cells.Item(47, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(47, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(47, 19).Value = 972
$ws.Cells.Item(47, 20).Value = 18
